{"js": "// Update the date stamp and the 25 multiplication problems to the new\n// day's values. Every old string is unique within the document, so a\n// targeted search-and-replace (matching whole word, case-sensitive) for\n// each pair is safe and keeps all other formatting untouched.\nconst replacements = [\n  [\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"],\n  [\"967\u00d72=\", \"838\u00d76=\"],\n  [\"516\u00d72=\", \"175\u00d74=\"],\n  [\"913\u00d75=\", \"531\u00d72=\"],\n  [\"632\u00d77=\", \"153\u00d75=\"],\n  [\"764\u00d76=\", \"728\u00d77=\"],\n  [\"746\u00d77=\", \"521\u00d79=\"],\n  [\"479\u00d72=\", \"804\u00d76=\"],\n  [\"780\u00d78=\", \"557\u00d74=\"],\n  [\"209\u00d74=\", \"965\u00d77=\"],\n  [\"475\u00d78=\", \"470\u00d79=\"],\n  [\"647\u00d79=\", \"365\u00d79=\"],\n  [\"102\u00d77=\", \"150\u00d77=\"],\n  [\"310\u00d73=\", \"559\u00d72=\"],\n  [\"801\u00d77=\", \"876\u00d76=\"],\n  [\"863\u00d72=\", \"340\u00d72=\"],\n  [\"941\u00d74=\", \"376\u00d78=\"],\n  [\"640\u00d76=\", \"846\u00d75=\"],\n  [\"396\u00d79=\", \"861\u00d73=\"],\n  [\"898\u00d77=\", \"227\u00d79=\"],\n  [\"734\u00d78=\", \"358\u00d75=\"],\n  [\"862\u00d74=\", \"559\u00d76=\"],\n  [\"288\u00d72=\", \"342\u00d73=\"],\n  [\"567\u00d75=\", \"855\u00d77=\"],\n  [\"143\u00d79=\", \"233\u00d74=\"],\n  [\"841\u00d79=\", \"159\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and the 25 multiplication problems to the new\n# day's values. Every old string is unique within the document, so a\n# targeted Find/Replace (match whole string, case-sensitive, no\n# wildcards) for each pair is safe and leaves all other formatting\n# untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"),\n    @(\"967\u00d72=\", \"838\u00d76=\"),\n    @(\"516\u00d72=\", \"175\u00d74=\"),\n    @(\"913\u00d75=\", \"531\u00d72=\"),\n    @(\"632\u00d77=\", \"153\u00d75=\"),\n    @(\"764\u00d76=\", \"728\u00d77=\"),\n    @(\"746\u00d77=\", \"521\u00d79=\"),\n    @(\"479\u00d72=\", \"804\u00d76=\"),\n    @(\"780\u00d78=\", \"557\u00d74=\"),\n    @(\"209\u00d74=\", \"965\u00d77=\"),\n    @(\"475\u00d78=\", \"470\u00d79=\"),\n    @(\"647\u00d79=\", \"365\u00d79=\"),\n    @(\"102\u00d77=\", \"150\u00d77=\"),\n    @(\"310\u00d73=\", \"559\u00d72=\"),\n    @(\"801\u00d77=\", \"876\u00d76=\"),\n    @(\"863\u00d72=\", \"340\u00d72=\"),\n    @(\"941\u00d74=\", \"376\u00d78=\"),\n    @(\"640\u00d76=\", \"846\u00d75=\"),\n    @(\"396\u00d79=\", \"861\u00d73=\"),\n    @(\"898\u00d77=\", \"227\u00d79=\"),\n    @(\"734\u00d78=\", \"358\u00d75=\"),\n    @(\"862\u00d74=\", \"559\u00d76=\"),\n    @(\"288\u00d72=\", \"342\u00d73=\"),\n    @(\"567\u00d75=\", \"855\u00d77=\"),\n    @(\"143\u00d79=\", \"233\u00d74=\"),\n    @(\"841\u00d79=\", \"159\u00d79=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Replacement.ClearFormatting()\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
